$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing rows (old rows 5, 6, 7) that no longer exist after the
# TPM data refresh - this also shifts the used range down to A1:T4.
$ws.Rows("5:7").Delete()

# Row 2: ECs -> FAPs (Agt/Lrp2) with refreshed TPM-derived values.
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07702833333333334
$ws.Range("H2").Value = 0.231085
$ws.Range("I2").Value = 0.02259036512642383
$ws.Range("J2").Value = 0.02259036512642383
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.0001512579705555556
$ws.Range("R2").Value = 0.001361321735
$ws.Range("S2").Value = 0.02259036512642383
$ws.Range("T2").Value = 0.02259036512642383

# Row 3: FAPs -> FAPs (Agt/Lrp2) with refreshed TPM-derived values.
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.276419666666667
$ws.Range("H3").Value = 9.829259
$ws.Range("I3").Value = 0.9608869019286738
$ws.Range("J3").Value = 0.9608869019286738
$ws.Range("M3").Value = 0.001963666666666667
$ws.Range("N3").Value = 0.005891
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.006433796085444446
$ws.Range("R3").Value = 0.057904164769
$ws.Range("S3").Value = 0.9608869019286738
$ws.Range("T3").Value = 0.9608869019286738

# Row 4: MuSCs -> FAPs (Agt/Lrp2) with refreshed TPM-derived values.
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.056339
$ws.Range("H4").Value = 0.169017
$ws.Range("I4").Value = 0.01652273294490242
$ws.Range("J4").Value = 0.01652273294490242
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.0001106310163333333
$ws.Range("R4").Value = 0.000995679147
$ws.Range("S4").Value = 0.01652273294490242
$ws.Range("T4").Value = 0.01652273294490242
